$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BU1").Value = "15-sep"

$ws.Range("BU2").Value = 0
$ws.Range("BU3").Value = 16.188640780434685
$ws.Range("BU4").Value = 18.368698306146527
$ws.Range("BU5").Value = 17.360829603157306
$ws.Range("BU6").Value = 0
$ws.Range("BU7").Value = 7.8096137130507701
$ws.Range("BU8").Value = 17.297045043200637
$ws.Range("BU9").Value = 10.830674401510734
$ws.Range("BU10").Value = 12.50323819001977
$ws.Range("BU11").Value = 13.801236067837319
$ws.Range("BU12").Value = 0
$ws.Range("BU13").Value = 10.42795331839403
$ws.Range("BU14").Value = 0
$ws.Range("BU15").Value = 0
$ws.Range("BU16").Value = 6.4314940272202277
$ws.Range("BU17").Value = 0
$ws.Range("BU18").Value = 0

$ws.Range("BV4").Select()
